# Add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计", cloned (layout/
#    formatting) from the existing "2022-Q3" sheet, then overwrite it with
#    the new quarter's numbers.
# 2. Insert a new top row into "总计" for the 2022-Q4 summary line, shifting
#    the older quarters down by one row (their data is unchanged).
# All other quarter sheets (2022-Q3 .. 2021-Q1) keep their original data and
# simply shift one tab to the right to make room for 2022-Q4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (keeps headers,
#    styles, number formats identical) and placing it right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fund A (010343)
$q4.Range("C2").Value = "华宝英国富时100指数A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.14"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.75"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "5.17"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0072"
$q4.Range("H2").Value = 4

# Fund C (010344)
$q4.Range("C3").Value = "华宝英国富时100指数C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.08"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "94.75"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "5.17"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0041"
$q4.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 2. Shift the "总计" rows down by one (bottom-up so nothing is clobbered)
#    and insert the 2022-Q4 summary line at the top of the table.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $srcRow = $total.Range("A" + $r + ":D" + $r)
    $dstRow = $total.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $srcRow.Copy($dstRow)
}

$total.Range("B2").Value = "2022-Q4"
# C2 (count) and D2 (value) for 2022-Q4 happen to match the values that were
# already sitting in row 2 before the shift (2, 0.01), so nothing else to do.

# ---------------------------------------------------------------------
# Restore the originally active tab (2021-Q1) as the selected sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
